$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Easy")

# Row 14: LeetCode #20 - Valid Parentheses
$ws.Range("A14").Value = 20
$ws.Range("B14").Value = "Valid Parentheses"
$ws.Range("C14").Value = "Stack & Queue"
$ws.Range("D14").Value = "Pedoe"
$ws.Range("E14").Value = "On-going"
$ws.Range("F14").Value = "Javascript"

# Row 15: LeetCode #503 - Next Greater Element II
$ws.Range("A15").Value = 503
$ws.Range("B15").Value = "Next Greater Element II"
$ws.Range("C15").Value = "Stack & Queue"
$ws.Range("D15").Value = "Pedoe"
$ws.Range("E15").Value = "On-going"
$ws.Range("F15").Value = "Javascript"

# Row 16: LeetCode #394 - Decoding String
$ws.Range("A16").Value = 394
$ws.Range("B16").Value = "Decoding String"
$ws.Range("C16").Value = "Stack & Queue"
$ws.Range("D16").Value = "Pedoe"
$ws.Range("E16").Value = "On-going"
$ws.Range("F16").Value = "Javascript"

# Formatting: columns A-D plain black Calibri 12
$ws.Range("A14:D16").Font.Color = 0

# Formatting: column E bold orange (status "On-going")
$rangeE = $ws.Range("E14:E16")
$rangeE.Font.Bold = $true
$rangeE.Font.Color = 3243501

# Formatting: column F bold blue (language "Javascript")
$rangeF = $ws.Range("F14:F16")
$rangeF.Font.Bold = $true
$rangeF.Font.Color = 12874308

$ws.Range("B18").Select() | Out-Null
